$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last refreshed" timestamp shown at the top of the sheet
$ws.Range("A1").Value = 'Datos actualizados a 29 de Marzo de 2020 a las 06:55'

# Refreshed country statistics. The underlying feed re-sorts all rows by
# "Casos totales" (column B) each time it is refreshed, so updated figures
# for a handful of countries shuffle the whole ranking. The table below is
# the resulting row-by-row state for every row whose contents move or change.
$rowData = @(
    @{Row=12; Country='Suiza'; Vals=@(14076, 0, 1595, 12217, 301, 0, 264)}
    @{Row=18; Country='Canada'; Vals=@(5655, 0, 508, 5087, 120, 0, 60)}
    @{Row=21; Country='Australia'; Vals=@(3969, 334, 170, 3783, 23, 2, 16)}
    @{Row=22; Country='Brasil'; Vals=@(3904, 0, 6, 3784, 296, 0, 114)}
    @{Row=101; Country='Sri Lanka'; Vals=@(115, 2, 9, 105, 5, 0, 1)}
    @{Row=102; Country='Honduras'; Vals=@(110, 15, 3, 106, 4, 0, 1)}
    @{Row=103; Country='Afganistan'; Vals=@(110, 0, 2, 104, 0, 0, 4)}
    @{Row=104; Country='Uzbekistan'; Vals=@(104, 0, 5, 97, 8, 0, 2)}
    @{Row=105; Country='Estado de Palestina'; Vals=@(104, 0, 18, 85, 0, 0, 1)}
    @{Row=106; Country='Camboya'; Vals=@(103, 4, 21, 82, 1, 0, 0)}
    @{Row=107; Country='Mauricio'; Vals=@(102, 0, 0, 100, 1, 0, 2)}
    @{Row=108; Country='Guadalupe'; Vals=@(102, 0, 17, 83, 4, 0, 2)}
    @{Row=109; Country='Costa de Marfil'; Vals=@(101, 0, 3, 98, 0, 0, 0)}
    @{Row=110; Country='Nigeria'; Vals=@(97, 0, 3, 93, 0, 0, 1)}
    @{Row=160; Country='Granada'; Vals=@(9, 2, 0, 9, 0, 0, 0)}
    @{Row=162; Country='Birmania'; Vals=@(8, 0, 0, 8, 0, 0, 0)}
    @{Row=163; Country='Seychelles'; Vals=@(8, 0, 0, 8, 0, 0, 0)}
    @{Row=164; Country='Surinam'; Vals=@(8, 0, 0, 8, 0, 0, 0)}
    @{Row=165; Country='Mozambique'; Vals=@(8, 0, 0, 8, 0, 0, 0)}
    @{Row=166; Country='Guinea'; Vals=@(8, 0, 0, 8, 0, 0, 0)}
    @{Row=167; Country='Islas Caimanes'; Vals=@(8, 0, 0, 7, 0, 0, 1)}
    @{Row=168; Country='Guyana'; Vals=@(8, 0, 0, 7, 0, 0, 1)}
    @{Row=169; Country='Namibia'; Vals=@(8, 0, 2, 6, 0, 0, 0)}
    @{Row=170; Country='Curazao'; Vals=@(8, 0, 2, 5, 0, 0, 1)}
    @{Row=171; Country='Antigua y Barbuda'; Vals=@(7, 0, 0, 7, 0, 0, 0)}
    @{Row=172; Country='Gabon'; Vals=@(7, 0, 0, 6, 0, 0, 1)}
    @{Row=173; Country='Zimbabue'; Vals=@(7, 0, 0, 6, 0, 0, 1)}
    @{Row=174; Country='Santa Sede'; Vals=@(6, 0, 0, 6, 0, 0, 0)}
    @{Row=175; Country='Eritrea'; Vals=@(6, 0, 0, 6, 0, 0, 0)}
    @{Row=176; Country='Benin'; Vals=@(6, 0, 0, 6, 0, 0, 0)}
    @{Row=186; Country='Butan'; Vals=@(4, 1, 0, 4, 0, 0, 0)}
    @{Row=187; Country='Islas Turcas y Caicos'; Vals=@(4, 0, 0, 4, 0, 0, 0)}
    @{Row=188; Country='Congo'; Vals=@(4, 0, 0, 4, 0, 0, 0)}
    @{Row=189; Country='Nicaragua'; Vals=@(4, 0, 0, 3, 0, 0, 1)}
    @{Row=190; Country='Republica de Africa Central'; Vals=@(3, 0, 0, 3, 0, 0, 0)}
)

foreach ($entry in $rowData) {
    $r = $entry.Row
    $ws.Cells.Item($r, 1).Value = $entry.Country
    for ($i = 0; $i -lt $entry.Vals.Length; $i++) {
        $ws.Cells.Item($r, 2 + $i).Value = $entry.Vals[$i]
    }
}
